$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 456/457, pushing the existing rows 456-553 down to 458-555
$ws.Rows("456:457").Insert()

# New row 456
$ws.Range("A456").Value = 4
$ws.Range("B456").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C456").Value = "Los Lagos"
$ws.Range("D456").Value = 44711
$ws.Range("E456").Value = 10
$ws.Range("F456").Value = 100112004
$ws.Range("G456").Value = "Cebolla"
$ws.Range("H456").Value = "Sin especificar"
$ws.Range("I456").Value = "1a (cosecha)"
$ws.Range("J456").Value = 250
$ws.Range("K456").Value = 9000
$ws.Range("L456").Value = 10000
$ws.Range("M456").Value = 9600
$ws.Range("N456").Value = "`$/malla 18 kilos"
$ws.Range("O456").Value = "Región de O'Higgins"
$ws.Range("P456").Value = 533
$ws.Range("Q456").Value = 18
$ws.Range("R456").Value = "Hortaliza"

# New row 457
$ws.Range("A457").Value = 4
$ws.Range("B457").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C457").Value = "Los Lagos"
$ws.Range("D457").Value = 44711
$ws.Range("E457").Value = 10
$ws.Range("F457").Value = 100112004
$ws.Range("G457").Value = "Cebolla"
$ws.Range("H457").Value = "Sin especificar"
$ws.Range("I457").Value = "1a (cosecha)"
$ws.Range("J457").Value = 250
$ws.Range("K457").Value = 12000
$ws.Range("L457").Value = 12000
$ws.Range("M457").Value = 12000
$ws.Range("N457").Value = "`$/malla 25 kilos"
$ws.Range("O457").Value = "Región del Maule"
$ws.Range("P457").Value = 480
$ws.Range("Q457").Value = 25
$ws.Range("R457").Value = "Hortaliza"
